$d = $word.ActiveDocument

# Locate the run of paragraphs to remove:
#   - the blank paragraph right after the bibliography's last entry
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - "© 2020 . Contact: ... Creative Commons Attribution"
# by scanning for their distinctive text, so the deletion is anchored to
# content rather than brittle hard-coded paragraph indices.
$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $startIdx = $i - 1
    }
    if ($t -like "*Creative Commons Attribution*") {
        $endIdx = $i
    }
}

if ($startIdx -gt 0 -and $endIdx -ge $startIdx) {
    $startPara = $d.Paragraphs.Item($startIdx)
    $endPara = $d.Paragraphs.Item($endIdx)
    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
